$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.209.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +13.49%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.674.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +8.26%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.30%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'309.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +9.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3743"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.23%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3448"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +7.72%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'47.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +16.43%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +6.19%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07309"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'20.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.69%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.73%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.59%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.677.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +8.45%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +5.38%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.17%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06731"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +9.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +12.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'16.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +7.92%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.149"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +7.09%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +5.09%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'24.165.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +13.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.70%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.359"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -9.33%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.660"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +16.03%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'151.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.60%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +9.76%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.862.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +8.40%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'126.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +6.91%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.438"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +22.41%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.9913"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +12.93%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.772"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +14.93%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.08509"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.66%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'12.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +16.53%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.06488"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +10.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.370"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +7.43%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'8.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +12.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.02360"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +10.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.279"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +4.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.2140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +6.99%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.6197"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +12.30%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.9986"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'13.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.67%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.810"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +6.35%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.5945"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +8.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'127.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.01%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.032"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +7.95%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.07168"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +8.36%  "
$ws.Range("E51").Style = "Normal"
